$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Gene Set" label text from "Combined" to "30-gene set" for all rows using it (E3:E7)
$ws.Range("E3:E7").Value = "30-gene set"

# Increase row height for data rows (2-7) from 20.1 to 24.95
$ws.Range("A2:A7").EntireRow.RowHeight = 24.95

# Update AUC and AUC_lower values for Fold1 (row 3)
$ws.Range("B3").Value = 0.99319727891156495
$ws.Range("C3").Value = 0.97707836111051805

# Update AUC and AUC_lower values for Fold2 (row 4)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1

# Update AUC and AUC_lower values for Fold3 (row 5)
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

# Update AUC and AUC_lower values for Fold4 (row 6)
$ws.Range("B6").Value = 0.97499999999999998
$ws.Range("C6").Value = 0.92290322306711503

# Update AUC and AUC_lower values for Fold5 (row 7)
$ws.Range("B7").Value = 0.96153846153846201
$ws.Range("C7").Value = 0.89702859416324399

$wb.Save()
